# B6-PowerPoint.pptx edit:
#  1. Re-colour the deck's theme (slide master / "Integral" design) so its
#     colour scheme becomes the stock Office colours (swap with the theme
#     that was previously only used by the Notes Master).
#  2. Re-apply the (built-in) table style {E93C81EF-2BED-463E-BC89-2B229A1D5B04}
#     to the three tables on slides 14, 15 and 16 (previously
#     {DB71E051-2E32-40C1-8C5B-2DAC43FF6F56}).

$p = $ppt.ActivePresentation

# --- 1. Theme re-colour ------------------------------------------------
# The presentation's theme colour scheme (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) is shared by every slide, so touching it once via the
# first slide updates the master theme for the whole deck.
$cs = $p.Slides.Item(1).ThemeColorScheme
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)
for ($i = 1; $i -le 12; $i++) {
    $cs.Colors($i).RGB = $officeColors[$i - 1]
}

# --- 2. Table style id swap ---------------------------------------------
$oldStyle = "{DB71E051-2E32-40C1-8C5B-2DAC43FF6F56}"
$newStyle = "{E93C81EF-2BED-463E-BC89-2B229A1D5B04}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyle)
        }
    }
}
